# This script updates the "Аркуш1" template table on the active worksheet.
# It (a) corrects the wording of two existing rows whose shared-string text
# changed, and (b) appends ten brand-new template rows (27-36) that were
# added to the sentence-generation table, re-using the existing row
# formatting (style) for the new cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Fix the wording of existing rows 15 and 18 (subject/verb/complement
#    rephrasing introduced by the edit).
# ---------------------------------------------------------------------
$ws.Range("A15").Value = "Студента"
$ws.Range("B15").Value = "запросили"
$ws.Range("C15").Value = "на інтернатуру"
$ws.Range("D15").Value = "в ?org_name"

$ws.Range("A18").Value = "Бухгалтер"
$ws.Range("B18").Value = "помилився"
$ws.Range("C18").Value = "в розрахунках"
$ws.Range("D18").Value = "заробітної плати працівників компанії ?org_name"

# ---------------------------------------------------------------------
# 2) Append the new template rows 27-36.  Copy the formatting of row 14
#    (an existing fully-populated A:E data row) down across the new
#    range first, so the new cells pick up the same cell style ("s=2")
#    as the rest of the table, then fill in the actual text values.
# ---------------------------------------------------------------------
$ws.Range("A14:E14").Copy()
$ws.Range("A27:E36").PasteSpecial(-4122)

$newRows = @(
    @{ Row = 27; A = "Він";       B = "жив";        C = "в офісі";              D = "компанії ?org_name";              E = "за адресою ?address" },
    @{ Row = 28; A = "Вона";      B = "привезла";    C = "документи";            D = "в офіс компанії ?org_name";       E = ", який знаходиться за адресою ?address" },
    @{ Row = 29; A = "Вони";      B = "збудували";   C = "офіс";                 D = " компанії ?org_name";             E = ", що знаходиться за адресою ?address" },
    @{ Row = 30; A = "Компанія";  B = "уклала";      C = "угоду";                D = "з ?org_name" },
    @{ Row = 31; A = "Таксист";   B = "підвіз";      C = "людину";               D = ", яка працює в ?org_name";        E = "за адресою ?address" },
    @{ Row = 32; A = "Банкір";    B = "надав";       C = "фінансові послуги";    D = "компанії ?org_name" },
    @{ Row = 33; A = "Військовий"; B = "отримав";    C = "фінансову допомогу";   D = "від компанії ?org_name" },
    @{ Row = 34; A = "?full_name"; B = "придбав";    C = "частину";              D = "компанії ?org_name";              E = ", яка зареєстрована за адресою ?address" },
    @{ Row = 35; A = "Літак";     B = "нажелить";    C = "компанії";             D = "?org_name" },
    @{ Row = 36; A = "Бізнесмен"; B = "продав";      C = "чатину";               D = "акцій компанії ?org_name" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    if ($r.ContainsKey("E")) {
        $ws.Range("E$row").Value = $r.E
    } else {
        $ws.Range("E$row").Clear()
    }
}
